$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("G4").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("F8").Value = 9522
$ws.Range("G8").Value = 0
$ws.Range("F10").Value = 33
$ws.Range("F12").Value = 1993
$ws.Range("F15").Value = 2633
$ws.Range("F16").Value = 130
$ws.Range("F17").Value = 3942
$ws.Range("F19").Value = 144
$ws.Range("F20").Value = 129
$ws.Range("F21").Value = 209
$ws.Range("F23").Value = 23
$ws.Range("F25").Value = 74
$ws.Range("F26").Value = 263
$ws.Range("F27").Value = 1068
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = 2374
$ws.Range("F30").Value = 1096
$ws.Range("F32").Value = 474
$ws.Range("F33").Value = 4321
$ws.Range("F35").Value = 199
$ws.Range("F36").Value = 352
$ws.Range("F37").Value = 181

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 194
$ws.Range("F3").Value = 977

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 194
$ws.Range("F4").Value = 977
$ws.Range("G6").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("F10").Value = 9522
$ws.Range("G10").Value = 0
$ws.Range("F12").Value = 33
$ws.Range("F14").Value = 1993
$ws.Range("F18").Value = 2633
$ws.Range("F19").Value = 130
$ws.Range("F20").Value = 3942
$ws.Range("F22").Value = 144
$ws.Range("F23").Value = 129
$ws.Range("F24").Value = 209
$ws.Range("F26").Value = 23
$ws.Range("F29").Value = 74
$ws.Range("F30").Value = 263
$ws.Range("F31").Value = 1068
$ws.Range("F32").Value = 3
$ws.Range("F33").Value = 2376
$ws.Range("F34").Value = 1096
$ws.Range("F36").Value = 474
$ws.Range("F37").Value = 4321
$ws.Range("F39").Value = 199
$ws.Range("F40").Value = 352
$ws.Range("F41").Value = 181
